$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 6667
$ws.Range("J32").Value = 7000.5
$ws.Range("L32").Value = 7000.5
$ws.Range("N32").Value = -7652.5
$ws.Range("H74").Value = 10000
$ws.Range("I74").Value = 10000
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 10000
$ws.Range("L74").Value = 0
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -9064
$ws.Range("H77").Value = 10000
$ws.Range("I77").Value = 10000
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 50000
$ws.Range("L77").Value = 0
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -45320
$ws.Range("H98").Value = 3482029.5
$ws.Range("I98").Value = 3789869
$ws.Range("J98").Value = 2004399.4
$ws.Range("K98").Value = 3789869
$ws.Range("L98").Value = 2004399.4
$ws.Range("M98").Value = -3788371
$ws.Range("N98").Value = -2007395.4
$ws.Range("H122").Value = 3482029.5
$ws.Range("I122").Value = 3789869
$ws.Range("J122").Value = 2004399.4
$ws.Range("K122").Value = 11369607
$ws.Range("L122").Value = 6013198.199999999
$ws.Range("M122").Value = -11367157
$ws.Range("N122").Value = -6018098.199999999
$ws.Range("H132").Value = 2917.6924
$ws.Range("I132").Value = 2527.6365
$ws.Range("J132").Value = 5063
$ws.Range("K132").Value = 7582.9095
$ws.Range("L132").Value = 15189
$ws.Range("M132").Value = -5052.9095
$ws.Range("N132").Value = -20249
$ws.Range("H137").Value = 35672.17
$ws.Range("I137").Value = 67049.17999999999
$ws.Range("K137").Value = 201147.54
$ws.Range("M137").Value = -198597.54

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1475.909
$ws.Range("I2").Value = 1497
$ws.Range("K2").Value = 1497
$ws.Range("M2").Value = -1384
$ws.Range("H5").Value = 198
$ws.Range("I5").Value = 198
$ws.Range("K5").Value = 198
$ws.Range("M5").Value = -86
$ws.Range("H30").Value = 985
$ws.Range("I30").Value = 985
$ws.Range("K30").Value = 985
$ws.Range("M30").Value = -835
$ws.Range("H61").Value = 4578.731
$ws.Range("I61").Value = 2707.75
$ws.Range("J61").Value = 10815.333
$ws.Range("K61").Value = 2707.75
$ws.Range("L61").Value = 10815.333
$ws.Range("M61").Value = -2495.75
$ws.Range("N61").Value = -11239.333
$ws.Range("H102").Value = 3327.889
$ws.Range("I102").Value = 3119
$ws.Range("J102").Value = 4999
$ws.Range("K102").Value = 3119
$ws.Range("L102").Value = 4999
$ws.Range("M102").Value = -1497
$ws.Range("N102").Value = -8243
$ws.Range("H116").Value = 1475.909
$ws.Range("I116").Value = 1497
$ws.Range("K116").Value = 1497
$ws.Range("M116").Value = 797
$ws.Range("H132").Value = 4718.5557
$ws.Range("I132").Value = 4718.5557
$ws.Range("K132").Value = 14155.6671
$ws.Range("M132").Value = -11625.6671
$ws.Range("H136").Value = 4578.731
$ws.Range("I136").Value = 2707.75
$ws.Range("J136").Value = 10815.333
$ws.Range("K136").Value = 8123.25
$ws.Range("L136").Value = 32445.999
$ws.Range("M136").Value = -5573.25
$ws.Range("N136").Value = -37545.999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1475.909
$ws.Range("I3").Value = 1497
$ws.Range("K3").Value = 1497
$ws.Range("M3").Value = -1383
$ws.Range("H4").Value = 198
$ws.Range("I4").Value = 198
$ws.Range("K4").Value = 198
$ws.Range("M4").Value = -83
$ws.Range("H86").Value = 2511
$ws.Range("I86").Value = 2454
$ws.Range("J86").Value = 2594.3076
$ws.Range("K86").Value = 2454
$ws.Range("L86").Value = 2594.3076
$ws.Range("M86").Value = -1331
$ws.Range("N86").Value = -4840.3076
$ws.Range("H89").Value = 2511
$ws.Range("I89").Value = 2454
$ws.Range("J89").Value = 2594.3076
$ws.Range("K89").Value = 12270
$ws.Range("L89").Value = 12971.538
$ws.Range("M89").Value = -6654
$ws.Range("N89").Value = -24203.538
$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("M96").ClearContents()
$ws.Range("H103").Value = 34391.4
$ws.Range("J103").Value = 34391.4
$ws.Range("L103").Value = 34391.4
$ws.Range("N103").Value = -36735.4
$ws.Range("H105").Value = 1437.3939
$ws.Range("J105").Value = 1085.2858
$ws.Range("L105").Value = 1085.2858
$ws.Range("N105").Value = -4579.2858
$ws.Range("H134").Value = 2287.7666
$ws.Range("I134").Value = 2263.2068
$ws.Range("K134").Value = 6789.6204
$ws.Range("M134").Value = -4254.6204
$ws.Range("H138").Value = 63347.5
$ws.Range("J138").Value = 63347.5
$ws.Range("L138").Value = 63347.5
$ws.Range("N138").Value = -73627.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H39").Value = 1074.25
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()
$ws.Range("H49").Value = 1074.25
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()
$ws.Range("H96").Value = 11783
$ws.Range("J96").Value = 11783
$ws.Range("L96").Value = 11783
$ws.Range("N96").Value = -17275
$ws.Range("H105").Value = 4872.4326
$ws.Range("I105").Value = 1426.75
$ws.Range("K105").Value = 1426.75
$ws.Range("M105").Value = 320.25
$ws.Range("H106").Value = 20555.5
$ws.Range("J106").Value = 20555.5
$ws.Range("L106").Value = 20555.5
$ws.Range("N106").Value = -23079.5
$ws.Range("H132").Value = 2065.4
$ws.Range("I132").Value = 2064.4546
$ws.Range("K132").Value = 6193.3638
$ws.Range("M132").Value = -3663.3638
$ws.Range("H134").Value = 4502.3413
$ws.Range("I134").Value = 4571.9487
$ws.Range("K134").Value = 13715.8461
$ws.Range("M134").Value = -11180.8461

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 4619.857
$ws.Range("I26").Value = 420
$ws.Range("J26").Value = 6299.8
$ws.Range("K26").Value = 1260
$ws.Range("L26").Value = 18899.4
$ws.Range("M26").Value = -972
$ws.Range("N26").Value = -19475.4
$ws.Range("H131").Value = 8773329
$ws.Range("I131").Value = 100000640
$ws.Range("J131").Value = 1471.7693
$ws.Range("K131").Value = 300001920
$ws.Range("L131").Value = 4415.3079
$ws.Range("M131").Value = -299996880
$ws.Range("N131").Value = -14495.3079
$ws.Range("H132").Value = 6064.846
$ws.Range("I132").Value = 6195.2915
$ws.Range("J132").Value = 4499.5
$ws.Range("K132").Value = 55757.6235
$ws.Range("L132").Value = 40495.5
$ws.Range("M132").Value = -53227.6235
$ws.Range("N132").Value = -45555.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 8683.286
$ws.Range("J20").Value = 9989.666999999999
$ws.Range("L20").Value = 9989.666999999999
$ws.Range("N20").Value = -10441.667
$ws.Range("H22").Value = 3966.6667
$ws.Range("I22").Value = 3900
$ws.Range("K22").Value = 3900
$ws.Range("M22").Value = -3605
$ws.Range("H27").Value = 3966.6667
$ws.Range("I27").Value = 3900
$ws.Range("K27").Value = 3900
$ws.Range("M27").Value = -3793
$ws.Range("H46").Value = 1832.0714
$ws.Range("I46").Value = 2063.0625
$ws.Range("J46").Value = 1689.9231
$ws.Range("K46").Value = 2063.0625
$ws.Range("L46").Value = 1689.9231
$ws.Range("M46").Value = -1875.0625
$ws.Range("N46").Value = -2065.9231
$ws.Range("H96").Value = 49985.668
$ws.Range("J96").Value = 49983
$ws.Range("L96").Value = 49983
$ws.Range("N96").Value = -55475
$ws.Range("H136").Value = 1299.4
$ws.Range("I136").Value = 1171.6666
$ws.Range("J136").Value = 1627.8572
$ws.Range("K136").Value = 3514.9998
$ws.Range("L136").Value = 4883.571599999999
$ws.Range("M136").Value = -964.9998000000001
$ws.Range("N136").Value = -9983.571599999999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 19714
$ws.Range("I14").Value = 16999.666
$ws.Range("J14").Value = 36000
$ws.Range("K14").Value = 16999.666
$ws.Range("L14").Value = 36000
$ws.Range("M14").Value = -16831.666
$ws.Range("N14").Value = -36336
$ws.Range("H132").Value = 1289.8182
$ws.Range("I132").Value = 1218.8
$ws.Range("K132").Value = 3656.4
$ws.Range("M132").Value = -1126.4
